$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C rows 2 through 103 hold a "Förändrad" (last-changed) date,
# stored as serial date 45180 (2023-09-11). Update every row to 45181
# (2023-09-12), matching the diff which bumps C2:C103 by exactly one day.
$lastRow = 103
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45181
}
